$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPIs")
$ws.Range("D6").Value = "1. Main Menu, 2. Features Menu, 3. Spirits List Menu, 4. Bottle List Menu, 5. After Dinner / Dessert Menu"
